# Minor clarifications to the "Interface NamedDecl" slide (Abstract Syntax
# Trees deck): expand the parenthetical that explains VarDecl / SingleVarDecl
# in the body placeholder's first paragraph.

$p = $ppt.ActivePresentation

# Locate the slide/shape containing the paragraph to edit (rather than
# trusting a hard-coded slide number) so the script is resilient.
$targetSlide = $null
$targetShape = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*Identifiers declared using*") {
                $targetSlide = $slide
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$para = $targetShape.TextFrame.TextRange.Paragraphs(1)

# Walk the runs of the paragraph and update the two plain-text runs that sit
# between the "VarDecl" / "SingleVarDecl" / "ParameterDecl" code-styled runs.
for ($k = 1; $k -le $para.Runs().Count; $k++) {
    $run = $para.Runs($k)
    if ($run.Text -eq " (") {
        $run.Text = " ( which we convert to a list of "
    }
    elseif ($run.Text -eq ") or ") {
        $run.Text = " as described later) or "
    }
}
